$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Pseudotime_1", "Pseudotime_2" and "Pseudotime_3" rows (rows 35-37) were
# removed from the table. Deleting these rows shifts every subsequent row
# (Brain_region:, Cortex, Lesion, nCount_RNA, ... 461) up by three, which also
# shrinks the sheet's used range from A1:D54 down to A1:D51.
$ws.Rows("35:37").Delete()
